# Inclusao de bola azul (logo) no kickoff
# Adds a blue oval ("Oval 3") autoshape to slide 1, matching the
# accent1-themed circle placed on the kickoff slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# msoShapeOval = 9. Left/Top/Width/Height are expressed in points;
# 8750105 / 3429000 / 2616591 / 2912013 EMU converted to points
# (1 pt = 12700 EMU).
$oval = $s.Shapes.AddShape(9, 688.9846456692913, 270, 206.0307874015748, 229.2923622047244)

# Give it the default "Colored Fill - Accent 1" look (solid accent1 fill,
# accent1 outline) so it renders as the intended blue ball.
$oval.Fill.Solid()
$oval.Fill.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1
$oval.Line.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1

# Centered text body / paragraph, matching the inserted shape's layout.
$oval.TextFrame.VerticalAnchor = 3          # msoAnchorCenter
$oval.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter
